# Add 2018 (AD) and 2019 (AE) data points for the ".1.0.0.0" population
# variable rows that were still missing them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AD4").Value = 391.2626588110001
$ws.Range("AE4").Value = 412.5551113083447

$ws.Range("AD5").Value = 394.2290174485
$ws.Range("AE5").Value = 416.7214275212921

$ws.Range("AD12").Value = 1.862144777000012
$ws.Range("AE12").Value = 3.139585926351273

$ws.Range("AD13").Value = 2.96635863749998
$ws.Range("AE13").Value = 4.166316212947322

$ws.Range("AD16").Value = 96.00819525389991
$ws.Range("AE16").Value = 99.4222724380695

$ws.Range("AD17").Value = 23.11901260229315
$ws.Range("AE17").Value = 24.24488483060908

$ws.Range("AD19").Value = 443.2633733300499
$ws.Range("AE19").Value = 457.6293185330783

$ws.Range("AD21").Value = 449.15827333005
$ws.Range("AE21").Value = 463.5242185330783

$ws.Range("AD23").Value = 843.38729077855
$ws.Range("AE23").Value = 880.2456460543705

$ws.Range("AD25").Value = 112.7745918852599
$ws.Range("AE25").Value = 117.015989963076
